# Avance muestreo de datos SpaOnline
$wb = $excel.ActiveWorkbook

# --- Sheet "Objetos de dominio": the "Sucursales" domain object description
#     now also mentions it holds the branches' *location* data.
$ws2 = $wb.Worksheets.Item("Objetos de dominio")
$ws2.Range("B2").Value = "Objeto de dominio que contiene la informacion de la ubicación de las sucursales del Spa"
$ws2.Rows.Item(2).RowHeight = 30

# --- Sheet "Sucursales": turn the empty 4-column template (Identificador /
#     Nombre / Fabricante / Combinacion unica) into a real sample of branch
#     location data: Identificador / Pais / Departamento / Ciudad /
#     Sucursal / Ubicacion / Combinacion unica.
$ws3 = $wb.Worksheets.Item("Sucursales")

# Column D used to be the last (highlighted) column holding the formula; now
# it becomes an ordinary middle data column like B/C, and the *new* column G
# becomes the highlighted formula column instead.
#
# Clone D's current ("highlighted formula column") formatting onto the new
# G column FIRST, then clone C's ("ordinary column") formatting onto D, so
# the styling moves along with the columns' new roles.
$ws3.Range("D1").Copy()
$ws3.Range("G1").PasteSpecial(-4122)
$ws3.Range("D2:D4").Copy()
$ws3.Range("G2:G4").PasteSpecial(-4122)

$ws3.Range("C1").Copy()
$ws3.Range("D1").PasteSpecial(-4122)
$ws3.Range("C2:C4").Copy()
$ws3.Range("D2:D4").PasteSpecial(-4122)

# Add the two other new columns (E/F), cloning formatting from column A.
$ws3.Range("A1").Copy()
$ws3.Range("E1:F1").PasteSpecial(-4122)
$ws3.Range("A2:A4").Copy()
$ws3.Range("E2:F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-purpose the existing header cells (B/C/D).
$ws3.Range("B1").Value = "País"
$ws3.Range("C1").Value = "Departamento"
$ws3.Range("D1").Value = "Ciudad"

$ws3.Range("E1").Value = "Sucursal"
$ws3.Range("F1").Value = "Ubicación"
$ws3.Range("G1").Value = "Combinación única"

# Fill in the sample data rows.
$ws3.Range("B2").Value = "Colombia"
$ws3.Range("C2").Value = "Antioquia"
$ws3.Range("D2").Value = "Medellín"
$ws3.Range("E2").Value = "El poblado"
$ws3.Range("F2").Value = "CL 10 43 A 29"

$ws3.Range("B3").Value = "Colombia"
$ws3.Range("C3").Value = "Cundinamarca"
$ws3.Range("D3").Value = "Bogotá"
$ws3.Range("E3").Value = "Chapinero"
$ws3.Range("F3").Value = "CL 63 9 36"

$ws3.Range("B4").Value = "Colombia"
$ws3.Range("C4").Value = "Atlántico"
$ws3.Range("D4").Value = "Cartagena"
$ws3.Range("E4").Value = "Bocagrande"
$ws3.Range("F4").Value = "CR 2 5 39"

# New "Combinacion unica" column concatenates all five data columns.
$ws3.Range("G2").Formula = "=+B2&""-""&C2&""-""&D2&""-""&E2&""-""&F2"
$ws3.Range("G3").Formula = "=+B3&""-""&C3&""-""&D3&""-""&E3&""-""&F3"
$ws3.Range("G4").Formula = "=+B4&""-""&C4&""-""&D4&""-""&E4&""-""&F4"

# Match the column widths Excel would have picked via "best fit" for the
# new text (engine rounds ColumnWidth to 1/6-character increments, so we
# solve for the input that lands closest to the real target width).
$ws3.Columns.Item(2).ColumnWidth = 9.0
$ws3.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws3.Columns.Item(4).ColumnWidth = 9.5
$ws3.Columns.Item(5).ColumnWidth = 10.666666666666666
$ws3.Columns.Item(6).ColumnWidth = 11.5
$ws3.Columns.Item(7).ColumnWidth = 50.666666666666664

# Make "Sucursales" the active sheet/tab and select the full sample range,
# matching the state Excel leaves behind after finishing the data entry.
$ws3.Activate()
$ws3.Range("A1:G4").Select()

# "Objetos de dominio" is no longer the tab on top.
$ws2.Range("B10").Select()
